$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c-demo_ui")

# Insert a new row above row 22 (logging.details.enabled) for the new
# uiActionLog.details.enabled config entry.
$ws.Rows.Item(22).Insert()

$ws.Cells.Item(22, 1).Value = "uiActionLog.details.enabled"
$ws.Cells.Item(22, 2).Value = "TRUE"

# Append new config row at the end of the sheet for the codeless plugin.
$ws.Cells.Item(27, 1).Value = "codeless.plugins"
$ws.Cells.Item(27, 2).Value = "com.tmobile.etp.codeless.plugin.reporting.CodelessPluginManager"
